$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 43
$ws1.Range("F3").Value = 316
$ws1.Range("F4").Value = 240
$ws1.Range("F5").Value = 2839
$ws1.Range("F6").Value = 1973
$ws1.Range("F8").Value = 132
$ws1.Range("F9").Value = 1035
$ws1.Range("F10").Value = 195
$ws1.Range("F11").Value = 239
$ws1.Range("F12").Value = 41

# Sheet "全部类型" (All Types) - same updates, plus row shifted by one
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 43
$ws4.Range("F3").Value = 316
$ws4.Range("F4").Value = 240
$ws4.Range("F5").Value = 2839
$ws4.Range("F6").Value = 1973
$ws4.Range("F9").Value = 132
$ws4.Range("F10").Value = 1035
$ws4.Range("F11").Value = 195
$ws4.Range("F12").Value = 239
$ws4.Range("F13").Value = 41
